# Add "Teacher Table" and "Student Table" blocks below the existing
# "Admin Table" block, each with a title row, a header row
# (id / name / email / password) and a PK/unique annotation row.
# Also backfills the PK/unique row for the pre-existing Admin Table and
# re-styles its header row to match the new look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restyle the existing Admin Table header row (row 4): still based on the
# "Calculation" cell style, but with the font/fill overridden to a bold
# dark-on-grey look instead of the default orange-on-light-grey.
function Style-HeaderRow($row) {
    $rng = $ws.Range("A" + $row + ":D" + $row)
    $rng.Style = "Calculation"
    $rng.Font.Bold = $true
    $rng.Font.ThemeColor = 1
    $rng.Interior.Color = 12566463
}

# Populate + style one "PK" / "unique" annotation row below a header row.
function Add-PkRow($row) {
    $rng = $ws.Range("A" + $row + ":D" + $row)
    $rng.Style = "Output"
    $ws.Cells.Item($row, 1).Value = "PK"
    $ws.Cells.Item($row, 3).Value = "unique"
    $ws.Range("A" + $row).Font.Bold = $false
    $ws.Range("C" + $row).Font.Bold = $false
}

# Build a full table block: title row (merged) + header row + PK row.
function Add-TableBlock($titleRow, $titleText) {
    $headerRow = $titleRow + 1
    $pkRow = $titleRow + 2

    $titleRange = $ws.Range("A" + $titleRow + ":D" + $titleRow)
    $titleRange.Merge() | Out-Null
    $ws.Cells.Item($titleRow, 1).Value = $titleText
    $titleRange.HorizontalAlignment = -4108
    $titleRange.Font.Bold = $true
    $titleRange.Font.Size = 14
    $titleRange.Font.ThemeColor = 1
    $ws.Rows.Item($titleRow).RowHeight = 18

    $ws.Cells.Item($headerRow, 1).Value = "id"
    $ws.Cells.Item($headerRow, 2).Value = "name"
    $ws.Cells.Item($headerRow, 3).Value = "email"
    $ws.Cells.Item($headerRow, 4).Value = "password"
    Style-HeaderRow $headerRow

    Add-PkRow $pkRow
}

# 1. Backfill the PK/unique row for the pre-existing Admin Table, and
#    restyle its header row.
Style-HeaderRow 4
Add-PkRow 5

# 2. Teacher Table block (rows 8-10).
Add-TableBlock 8 "Teacher Table"

# 3. Student Table block (rows 13-15).
Add-TableBlock 13 "Student Table"

$ws.Range("D26").Select() | Out-Null
